# Apply edits to match target diff: remove duplicate User_TC003 row,
# refresh description/result wording, and fix a couple of run-mode values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 in the original sheet duplicates row 4 (User_TC003 accept-invite case).
# Delete it so everything below shifts up by one row.
$ws.Rows(5).Delete()

$b2 = @'
To validate invite via email based on following conditions.
1) Launch application and login as global admin.
2) Click on Users Menu.
3) Click on Add User button.
4) Enter valid details in required field (Full Name, Email).
5) Select any one role (Carrier, Shipper Admin, Shipper User, Driver (Full) and Driver (Limited) .
6) Select any one value from Carrier Drop Down field.
7) Click on Invite button and check invitation email sent successfully.
8.) Clickon Resend Invitation and Check if invitation sent again.
'@

$ws.Range("B2").Value = $b2

$b3 = @'
To validate invite via email based on following conditions
1) Launch application and login as global admin.
2) Click on Users Menu.
3) Click on Add User button.
4) Enter valid details in required field (Full Name, Email)
5) Select any one role (Carrier, Shipper Admin, Shipper User, Driver (Full) and Driver (Limited) .
6) Select any one value from Carrier Drop Down field
7) Click on Invite button and check invitation email sent successfully.
8.) Clickon cancel invitation and check if invitation cancelled.
'@

$ws.Range("B3").Value = $b3

$b4 = @'
To validate invite via email for Shipper Admin  based on following conditions
1) Launch application and login as global admin.
2) Click on Users Menu.
3) Click on Add User button.
4) Enter valid details in required field (Full Name, Email)
5) Select any one role (Carrier, Shipper Admin, Shipper User, Driver (Full) and Driver (Limited) .
6) Select any one value from Carrier Drop Down field.
7) Click on Invite button and check invitation email sent successfully.
8.) Accet email invitation and check if invitation acceted.
'@

$ws.Range("B4").Value = $b4

$b5 = @'
Validate Forgot password based on following conditions:
1. Launch application and Set UserName.
2. Clickon Forgot Password.
3. Verify email and reset password.
4. Login with new password.
5. Check whether user is able to login with new password.
'@

$ws.Range("B5").Value = $b5

$b6 = @'
Validate delete any  user based on following conditions:
1. Launch application .
2. Login as global admin.
3. Go to Users and search accepted user(except global admin)
4. Clickon delete.
5. Check whether deleted user exist.
'@

$ws.Range("B6").Value = $b6

$b7 = @'
Validate whehter Global admin is able to add new load using Shipper platform on following conditions.
a.) Launch and login applcation as Global admin
b.) Click on add new load button.
C.) Set Carrier Name, Load Date, Shipper, Rate, Rate UOM and Commodity.
d.) click on save and search for record in AG grid.
e.) click on Edit and set Shipper contact and click on save.
f.) Search for edited record in AG grid
g.) Click on delete button.
h.) Check whether record exist in AG grid after delete.
'@

$ws.Range("B7").Value = $b7

$b8 = @'
Validate whehter carrier user can add Scoular loads for payment using full submit
1) Enter valid user id and Password and click Login button.
2) Click on Add New Load button from Load menu.
3) Enter valid details in all required field and click Save button
4) Now loads are saved successfully.
5) Upload an Origin and Destination ticket image or PDF document for corresponding load.
6) Observe Ready to Submit Load icon in grid should change to green color.
7) Click on Submit Load button.
8) Select any option and click Submit button.
'@

$ws.Range("B8").Value = $b8

# Loads_TC001 (row 7): "Webtable validated successfully" -> "Load added successfully"
$ws.Range("D7").Value = "Load added successfully"

# Loads_TC002 (row 8) now documents the Scoular full-submit case, run as "Yes"
$ws.Range("C8").Value = "Yes"
$ws.Range("D8").Value = "Load submitted successfully"

# Loads_TC003 (row 9) run mode flips from YES to NO
$ws.Range("C9").Value = "NO"

# Match the author's final selection (Run Mode cell for Loads_TC003)
$ws.Range("C9").Select()
